# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "67.064.13"
$ws.Range("E2").Value = "  +1.37%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "3.116.20"
$ws.Range("E3").Value = "  +2.98%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.12%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.89"
$ws.Range("E5").Value = "  +0.72%  "

# --- Row 6 (Solana) ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.09"
$ws.Range("E6").Value = "  +2.79%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  -0.10%  "

# --- Row 8 (LidoStakedEther) ---
$ws.Range("D8").Value = "3.114.93"
$ws.Range("E8").Value = "  +3.11%  "

# --- Row 9 (XRP) ---
$ws.Range("E9").Value = "  +0.73%  "

# --- Row 10 (Toncoin) ---
$ws.Range("E10").Value = "  -3.01%  "

# --- Row 11 (Dogecoin) ---
$ws.Range("E11").Value = "  +1.54%  "

# --- Row 12 (Cardano) ---
$ws.Range("E12").Value = "  +0.13%  "

# --- Row 13 (ShibaInu) ---
$ws.Range("E13").Value = "  +1.07%  "

# --- Row 14 (Avalanche) ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.21"
$ws.Range("E14").Value = "  +1.89%  "

# --- Row 15 (TRON) ---
$ws.Range("E15").Value = "  +0.20%  "

# --- Row 16 (WrappedliquidstakedEther2.0) ---
$ws.Range("D16").Value = "3.633.56"
$ws.Range("E16").Value = "  +3.04%  "

# --- Row 17 (WrappedBTC) ---
$ws.Range("D17").Value = "67.038.04"
$ws.Range("E17").Value = "  +1.36%  "

# --- Row 18 (Polkadot) ---
$ws.Range("E18").Value = "  -0.87%  "

# --- Row 19 (WrappedEther) ---
$ws.Range("D19").Value = "3.116.28"
$ws.Range("E19").Value = "  +2.95%  "

# --- Row 20 (Chainlink) ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.21"
$ws.Range("E20").Value = "  -1.32%  "

# --- Row 21 (BitcoinCash) ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.19"
$ws.Range("E21").Value = "  +3.60%  "

# --- Row 22 (Polygon) ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +1.60%  "

# --- Row 23 (Uniswap) ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  +2.00%  "

# --- Row 24 (Litecoin) ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.48"
$ws.Range("E24").Value = "  +1.46%  "

# --- Row 25 (InternetComputer(DFINITY)) ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.41"
$ws.Range("E25").Value = "  +4.33%  "

# --- Row 26 (Fetch.AI) ---
$ws.Range("E26").Value = "  +4.33%  "

# --- Row 27 (RenderToken) ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  +0.59%  "

# --- Row 28 (Dai) ---
$ws.Range("E28").Value = "  +0.08%  "

# --- Row 29 (NEARProtocol) ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  -2.39%  "

# --- Row 30 (ImmutableX) ---
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.43"
$ws.Range("E30").Value = "  -0.69%  "

# --- Row 31 (PancakeSwap) ---
$ws.Range("E31").Value = "  +2.72%  "

# --- Row 32 (EthereumClassic) ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.94"
$ws.Range("E32").Value = "  +3.02%  "

# --- Row 33 (PEPE) ---
$ws.Range("E33").Value = "  +1.16%  "

# --- Row 34 (Hedera) ---
$ws.Range("E34").Value = "  -4.40%  "

# --- Row 35 (FirstDigitalUSD) ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.01%  "

# --- Row 36 (Filecoin) ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  +1.47%  "

# --- Row 37 (Mantle) ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.990"
$ws.Range("E37").Value = "  +0.27%  "

# --- Row 38 (Arweave) ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.57"
$ws.Range("E38").Value = "  -0.94%  "

# --- Row 39 (Stacks) ---
$ws.Range("E39").Value = "  +2.99%  "

# --- Row 40 (OKB) ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.15"
$ws.Range("E40").Value = "  +1.15%  "

# --- Row 41 (TheGraph) ---
$ws.Range("E41").Value = "  +2.43%  "

# --- Row 42 (Kaspa) ---
$ws.Range("E42").Value = "  +1.30%  "

# --- Row 43 (Cosmos) ---
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.67"
$ws.Range("E43").Value = "  +0.66%  "

# --- Row 44 (dogwifhat) ---
$ws.Range("E44").Value = "  -1.19%  "

# --- Row 45 (Maker) ---
$ws.Range("D45").Value = "2.844.08"
$ws.Range("E45").Value = "  +4.48%  "

# --- Rows 46/47 swap: Bittensor/VeChain trade places (with updated values) ---
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0359"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "385.51"
$ws.Range("E47").Value = "  +1.04%  "

# --- Row 48 (Monero) ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "137.07"
$ws.Range("E48").Value = "  +2.01%  "

# --- Row 49 (USDe) ---
$ws.Range("E49").Value = "  +0.00%  "

# --- Row 50 (InjectiveProtocol) ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.13"
$ws.Range("E50").Value = "  +1.86%  "

# --- Row 51 (ThetaToken) ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  +0.49%  "
